# Auto-applied market-price / profit refresh for Kujata_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29186.223
$ws.Range("I21").Value = 27666.666
$ws.Range("J21").Value = 29946
$ws.Range("K21").Value = 27666.666
$ws.Range("L21").Value = 29946
$ws.Range("M21").Value = -27198.666
$ws.Range("N21").Value = -30882
$ws.Range("H23").Value = 29186.223
$ws.Range("I23").Value = 27666.666
$ws.Range("J23").Value = 29946
$ws.Range("K23").Value = 27666.666
$ws.Range("L23").Value = 29946
$ws.Range("M23").Value = -27432.666
$ws.Range("N23").Value = -30414
$ws.Range("H28").Value = 4344.9375
$ws.Range("I28").Value = 4027.0833
$ws.Range("K28").Value = 4027.0833
$ws.Range("M28").Value = -3542.0833
$ws.Range("H33").Value = 296.07144
$ws.Range("I33").Value = 195.36363
$ws.Range("K33").Value = 195.36363
$ws.Range("M33").Value = 33.63637
$ws.Range("H113").Value = 18184244
$ws.Range("I113").Value = 22224522
$ws.Range("J113").Value = 2999.5
$ws.Range("K113").Value = 22224522
$ws.Range("L113").Value = 2999.5
$ws.Range("M113").Value = -22221268
$ws.Range("N113").Value = -9507.5
$ws.Range("H132").Value = 8341670
$ws.Range("I132").Value = 10422506
$ws.Range("J132").Value = 18326.5
$ws.Range("K132").Value = 31267518
$ws.Range("L132").Value = 54979.5
$ws.Range("M132").Value = -31264988
$ws.Range("N132").Value = -60039.5
$ws.Range("H137").Value = 2402.3062
$ws.Range("I137").Value = 1994.8518
$ws.Range("J137").Value = 2902.3635
$ws.Range("K137").Value = 5984.555399999999
$ws.Range("L137").Value = 8707.0905
$ws.Range("M137").Value = -3434.555399999999
$ws.Range("N137").Value = -13807.0905
$ws.Range("H138").Value = 2738.8406
$ws.Range("J138").Value = 2641.2307
$ws.Range("L138").Value = 7923.6921
$ws.Range("N138").Value = -18203.6921
$ws.Range("H140").Value = 37415
$ws.Range("J140").Value = 37415
$ws.Range("L140").Value = 37415
$ws.Range("N140").Value = -47775
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1500
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1387
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 11208.671
$ws.Range("I32").Value = 8532.161
$ws.Range("J32").Value = 19505.85
$ws.Range("K32").Value = 8532.161
$ws.Range("L32").Value = 19505.85
$ws.Range("M32").Value = -8245.161
$ws.Range("N32").Value = -20079.85
$ws.Range("H45").Value = 1318.3334
$ws.Range("I45").Value = 903.3333
$ws.Range("J45").Value = 1733.3334
$ws.Range("K45").Value = 903.3333
$ws.Range("L45").Value = 1733.3334
$ws.Range("M45").Value = -526.3333
$ws.Range("N45").Value = -2487.3334
$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 794
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 4358.9585
$ws.Range("I132").Value = 4249.3125
$ws.Range("J132").Value = 4578.25
$ws.Range("K132").Value = 12747.9375
$ws.Range("L132").Value = 13734.75
$ws.Range("M132").Value = -10217.9375
$ws.Range("N132").Value = -18794.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1500
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1386
$ws.Range("N3").ClearContents()
$ws.Range("H99").Value = 90910290
$ws.Range("J99").Value = 1474.5
$ws.Range("L99").Value = 1474.5
$ws.Range("N99").Value = -4470.5
$ws.Range("H107").Value = 976.7143
$ws.Range("I107").Value = 905.5454999999999
$ws.Range("J107").Value = 1237.6666
$ws.Range("K107").Value = 905.5454999999999
$ws.Range("L107").Value = 1237.6666
$ws.Range("M107").Value = 1014.4545
$ws.Range("N107").Value = -5077.6666
$ws.Range("H134").Value = 5447.5654
$ws.Range("I134").Value = 1140.2106
$ws.Range("J134").Value = 25907.5
$ws.Range("K134").Value = 3420.6318
$ws.Range("L134").Value = 77722.5
$ws.Range("M134").Value = -885.6318000000001
$ws.Range("N134").Value = -82792.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125001550
$ws.Range("I16").Value = 166668160
$ws.Range("J16").Value = 1700
$ws.Range("K16").Value = 166668160
$ws.Range("L16").Value = 1700
$ws.Range("M16").Value = -166667873
$ws.Range("N16").Value = -2274
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H107").Value = 734.5599999999999
$ws.Range("J107").Value = 1858.1666
$ws.Range("L107").Value = 1858.1666
$ws.Range("N107").Value = -5698.1666
$ws.Range("H113").Value = 125001550
$ws.Range("I113").Value = 166668160
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 166668160
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = -166665990
$ws.Range("N113").Value = -6040
$ws.Range("H132").Value = 1534.5333
$ws.Range("I132").Value = 1164.5294
$ws.Range("J132").Value = 2018.3846
$ws.Range("K132").Value = 3493.5882
$ws.Range("L132").Value = 6055.1538
$ws.Range("M132").Value = -963.5881999999997
$ws.Range("N132").Value = -11115.1538
$ws.Range("H135").Value = 37070.832
$ws.Range("J135").Value = 37070.832
$ws.Range("L135").Value = 37070.832
$ws.Range("N135").Value = -47210.832
$ws.Range("H138").Value = 129693.336
$ws.Range("J138").Value = 129693.336
$ws.Range("L138").Value = 129693.336
$ws.Range("N138").Value = -139973.336
$ws.Range("H139").Value = 44390
$ws.Range("J139").Value = 44390
$ws.Range("L139").Value = 44390
$ws.Range("N139").Value = -54670
$ws.Range("H140").Value = 55266.668
$ws.Range("J140").Value = 55266.668
$ws.Range("L140").Value = 55266.668
$ws.Range("N140").Value = -65626.66800000001
$ws.Range("H141").Value = 283659.34
$ws.Range("J141").Value = 296462.03
$ws.Range("L141").Value = 296462.03
$ws.Range("N141").Value = -306822.03
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 172348.39
$ws.Range("H107").Value = 9707.091
$ws.Range("J107").Value = 14914
$ws.Range("L107").Value = 44742
$ws.Range("N107").Value = -48582
$ws.Range("H137").Value = 26792076
$ws.Range("J137").Value = 11302.692
$ws.Range("L137").Value = 33908.076
$ws.Range("N137").Value = -44108.076
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4281.8887
$ws.Range("I102").Value = 3152.75
$ws.Range("J102").Value = 5185.2
$ws.Range("K102").Value = 3152.75
$ws.Range("L102").Value = 5185.2
$ws.Range("M102").Value = -1530.75
$ws.Range("N102").Value = -8429.200000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 23170
$ws.Range("J64").Value = 23170
$ws.Range("L64").Value = 23170
$ws.Range("N64").Value = -23620
$ws.Range("H67").Value = 23170
$ws.Range("J67").Value = 23170
$ws.Range("L67").Value = 23170
$ws.Range("N67").Value = -24730
$ws.Range("H115").Value = 29000
$ws.Range("J115").Value = 29000
$ws.Range("L115").Value = 29000
$ws.Range("N115").Value = -31350
$ws.Range("H132").Value = 2556.7856
$ws.Range("I132").Value = 2207.6667
$ws.Range("K132").Value = 6623.000100000001
$ws.Range("M132").Value = -4093.000100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 30007.6
$ws.Range("J28").Value = 30007.6
$ws.Range("L28").Value = 30007.6
$ws.Range("N28").Value = -30703.6
$ws.Range("H63").Value = 15018.777
$ws.Range("J63").Value = 15867.875
$ws.Range("L63").Value = 15867.875
$ws.Range("N63").Value = -17115.875
$ws.Range("H66").Value = 15018.777
$ws.Range("J66").Value = 15867.875
$ws.Range("L66").Value = 47603.625
$ws.Range("N66").Value = -53843.625
$ws.Range("H136").Value = 1116.2222
$ws.Range("I136").Value = 739.4666999999999
$ws.Range("K136").Value = 2218.4001
$ws.Range("M136").Value = 331.5999000000002
